# changed copper lantern textures
#
# Adds a new changelog entry ("Changed copper lantern textures") below the
# existing "Dropped OptiFine requirements ..." row, and moves the sheet's
# active selection to the next empty row (A12) as left by the author after
# typing the new line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New changelog line -> row 11, column A (same style as the rows above it,
# which Excel inherits automatically from column A's formatting).
$ws.Range("A11").Value = "Changed copper lantern textures"

# Leave the selection where the author's cursor ended up after adding the
# new row.
$ws.Range("A12").Select() | Out-Null
